# Insert a new data row at row 15 (shifting existing rows 15-32 down to 16-33)
# and populate it with a new weekly price record for "Haba" (Vega Monumental
# Concepcion), matching the rest of the constant columns already present in
# the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 15..32 down by one row, preserving formatting (Excel default
# behaviour copies the format of the row above into the newly inserted row).
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new record's data. The
# non-varying columns (A, B, C, E, F, G, H, I, N, Q, R) use the same values
# as every other row in this dataset.
$ws.Range("A15").Value2 = 11
$ws.Range("B15").Value2 = "Vega Monumental Concepción"
$ws.Range("C15").Value2 = "Bíobío"
$ws.Range("D15").Value2 = 44799
$ws.Range("E15").Value2 = 8
$ws.Range("F15").Value2 = 100112026
$ws.Range("G15").Value2 = "Haba"
$ws.Range("H15").Value2 = "Sin especificar"
$ws.Range("I15").Value2 = "Primera"
$ws.Range("J15").Value2 = 100
$ws.Range("K15").Value2 = 10000
$ws.Range("L15").Value2 = 11000
$ws.Range("M15").Value2 = 10500
$ws.Range("N15").Value2 = "$/saco 25 kilos"
$ws.Range("O15").Value2 = "Región de Coquimbo"
$ws.Range("P15").Value2 = 420
$ws.Range("Q15").Value2 = 25
$ws.Range("R15").Value2 = "Hortaliza"
